$wb = $excel.ActiveWorkbook

# CreateUser: ResultProd Pass, DateProd updated, ResultDemo stays Pass
$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed May 21 14:17:20 IST 2025"
$ws.Range("C2").Value = "Pass"

# FindUser: ResultProd Pass, DateProd updated, ResultDemo Fail
$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed May 21 14:17:54 IST 2025"
$ws.Range("C2").Value = "Fail"

# ModifyUser: ResultProd Pass, DateProd updated, ResultDemo Pass
$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed May 21 14:18:25 IST 2025"
$ws.Range("C2").Value = "Pass"

# ModifyUserPwd: ResultProd Pass, DateProd updated, ResultDemo Fail
$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed May 21 14:19:00 IST 2025"
$ws.Range("C2").Value = "Fail"

# AddDeleteRole: DateProd updated only
$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Wed May 21 14:15:11 IST 2025"

# SearchRole: DateProd updated only
$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Wed May 21 14:15:57 IST 2025"

$wb.Save()
